$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.873.29"
$ws.Range("E2").Value = "  -1.53%  "
$ws.Range("D3").Value = "2.358.41"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'325.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").Value = "'103.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.624"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'40.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.57%  "
$ws.Range("D11").Value = "'0.0925"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.64%  "
$ws.Range("E12").Value = "  -2.26%  "
$ws.Range("E13").Value = "  -3.42%  "
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("D15").Value = "'16.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.09%  "
$ws.Range("D16").Value = "2.710.99"
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("D17").Value = "2.366.09"
$ws.Range("E17").Value = "  -2.01%  "
$ws.Range("D18").Value = "42.719.26"
$ws.Range("E18").Value = "  -1.91%  "
$ws.Range("E19").Value = "  +8.96%  "
$ws.Range("E20").Value = "  -2.20%  "
$ws.Range("D21").Value = "'77.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.98%  "
$ws.Range("E22").Value = "  +3.72%  "
$ws.Range("D23").Value = "'266.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.71%  "
$ws.Range("D24").Value = "'2.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.02%  "
$ws.Range("D25").Value = "'9.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.61%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").Value = "'11.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.02%  "
$ws.Range("D28").Value = "'23.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.49%  "
$ws.Range("E29").Value = "  -1.20%  "
$ws.Range("D30").Value = "'174.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("E31").Value = "  -2.16%  "
$ws.Range("D32").Value = "'6.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.20%  "
$ws.Range("D33").Value = "'0.0902"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.57%  "
$ws.Range("D34").Value = "'35.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -9.12%  "
$ws.Range("D35").Value = "'0.134"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.78%  "
$ws.Range("E36").Value = "  +6.58%  "
$ws.Range("E37").Value = "  -7.95%  "
$ws.Range("E38").Value = "  -3.25%  "
$ws.Range("D39").Value = "'3.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.24%  "
$ws.Range("D40").Value = "'2.72"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.26%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'1.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.19%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.237"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.34%  "
$ws.Range("D43").Value = "'70.80"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.54%  "
$ws.Range("D44").Value = "'94.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +24.39%  "
$ws.Range("D45").Value = "'120.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.90%  "
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("D47").Value = "'11.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.50%  "
$ws.Range("D48").Value = "'5.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("D49").Value = "'9.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.14%  "
$ws.Range("E50").Value = "  -3.35%  "
$ws.Range("D51").Value = "'0.100"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.09%  "
